$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.163240432739258
$ws.Range("B1").Value = 2.421838521957397
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.376255750656128
$ws.Range("E1").Value = 1.234051942825317
